$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("numero", "lat", "lon", "data_arribada", "inici_simptomes", "pais_visitat", "estat", "centre_dia", "any", "prov", "edatany", "observacio")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Range("L2").Select()
